# Insert a new "Version" column in front of the existing Code/Description/
# Definition columns. Excel shifts the existing A/B/C columns to B/C/D
# (including the trailing blank cell in what is now D12).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(1).Insert()

# New header
$ws.Range("A1").Value = "Version"

# Every data row gets a "1.0" version tag. Writing the literal string "1.0"
# through .Value lets Excel's General-format numeric inference turn it into
# the number 1, so instead build it as a text formula and convert the
# formulas to static values (Copy + PasteSpecial values) — this keeps the
# values as real text ("1.0") stored as shared strings, matching how the
# sheet already stores every other code/description value, without forcing
# a new (Text) cell style onto the range.
$ws.Range("A2:A12").Formula = "=""1.0"""
$ws.Range("A2:A12").Copy()
$ws.Range("A2:A12").PasteSpecial(-4163)
